# Update the NATMI ligand-receptor pair results (Nxph3-Nrxn1) with the
# newly recomputed TPM-based statistics, and add the new "Resolving-Mac"
# target cluster rows for both sending clusters (ECs, MuSCs).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Nxph3"
$ws.Range("C2").Value = "Nrxn1"
$ws.Range("D2").Value = "FAPs"
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.1275786666666667
$ws.Range("H2").Value = 0.382736
$ws.Range("I2").Value = 0.2993455218931061
$ws.Range("J2").Value = 0.2993455218931061
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.2191816666666667
$ws.Range("N2").Value = 0.657545
$ws.Range("O2").Value = 0.8459226744507667
$ws.Range("P2").Value = 0.8459226744507669
$ws.Range("Q2").Value = 0.02796290479111111
$ws.Range("R2").Value = 0.25166614312
$ws.Range("S2").Value = 0.2532231644646769
$ws.Range("T2").Value = 0.2532231644646769

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Nxph3"
$ws.Range("C3").Value = "Nrxn1"
$ws.Range("D3").Value = "MuSCs"
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.1275786666666667
$ws.Range("H3").Value = 0.382736
$ws.Range("I3").Value = 0.2993455218931061
$ws.Range("J3").Value = 0.2993455218931061
$ws.Range("K3").Value = 2
$ws.Range("L3").Value = 0.6666666666666666
$ws.Range("M3").Value = 0.03648100000000001
$ws.Range("N3").Value = 0.109443
$ws.Range("O3").Value = 0.1407969268413801
$ws.Range("P3").Value = 0.1407969268413801
$ws.Range("Q3").Value = 0.004654197338666668
$ws.Range("R3").Value = 0.04188777604800001
$ws.Range("S3").Value = 0.0421469295462784
$ws.Range("T3").Value = 0.04214692954627841

# Row 4
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Nxph3"
$ws.Range("C4").Value = "Nrxn1"
$ws.Range("D4").Value = "Resolving-Mac"
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.1275786666666667
$ws.Range("H4").Value = 0.382736
$ws.Range("I4").Value = 0.2993455218931061
$ws.Range("J4").Value = 0.2993455218931061
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 0.6666666666666666
$ws.Range("M4").Value = 0.003441
$ws.Range("N4").Value = 0.010323
$ws.Range("O4").Value = 0.0132803987078531
$ws.Range("P4").Value = 0.0132803987078531
$ws.Range("Q4").Value = 0.000438998192
$ws.Range("R4").Value = 0.003950983728
$ws.Range("S4").Value = 0.003975427882150817
$ws.Range("T4").Value = 0.003975427882150818

# Row 5
$ws.Range("A5").Value = "MuSCs"
$ws.Range("B5").Value = "Nxph3"
$ws.Range("C5").Value = "Nrxn1"
$ws.Range("D5").Value = "FAPs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 0.2986133333333333
$ws.Range("H5").Value = 0.89584
$ws.Range("I5").Value = 0.700654478106894
$ws.Range("J5").Value = 0.700654478106894
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.2191816666666667
$ws.Range("N5").Value = 0.657545
$ws.Range("O5").Value = 0.8459226744507667
$ws.Range("P5").Value = 0.8459226744507669
$ws.Range("Q5").Value = 0.0654505680888889
$ws.Range("R5").Value = 0.5890551128
$ws.Range("S5").Value = 0.5926995099860899
$ws.Range("T5").Value = 0.59269950998609

# Row 6
$ws.Range("A6").Value = "MuSCs"
$ws.Range("B6").Value = "Nxph3"
$ws.Range("C6").Value = "Nrxn1"
$ws.Range("D6").Value = "MuSCs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 0.2986133333333333
$ws.Range("H6").Value = 0.89584
$ws.Range("I6").Value = 0.700654478106894
$ws.Range("J6").Value = 0.700654478106894
$ws.Range("K6").Value = 2
$ws.Range("L6").Value = 0.6666666666666666
$ws.Range("M6").Value = 0.03648100000000001
$ws.Range("N6").Value = 0.109443
$ws.Range("O6").Value = 0.1407969268413801
$ws.Range("P6").Value = 0.1407969268413801
$ws.Range("Q6").Value = 0.01089371301333334
$ws.Range("R6").Value = 0.09804341712
$ws.Range("S6").Value = 0.0986499972951017
$ws.Range("T6").Value = 0.09864999729510172

# Row 7
$ws.Range("A7").Value = "MuSCs"
$ws.Range("B7").Value = "Nxph3"
$ws.Range("C7").Value = "Nrxn1"
$ws.Range("D7").Value = "Resolving-Mac"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 0.2986133333333333
$ws.Range("H7").Value = 0.89584
$ws.Range("I7").Value = 0.700654478106894
$ws.Range("J7").Value = 0.700654478106894
$ws.Range("K7").Value = 2
$ws.Range("L7").Value = 0.6666666666666666
$ws.Range("M7").Value = 0.003441
$ws.Range("N7").Value = 0.010323
$ws.Range("O7").Value = 0.0132803987078531
$ws.Range("P7").Value = 0.0132803987078531
$ws.Range("Q7").Value = 0.00102752848
$ws.Range("R7").Value = 0.00924775632
$ws.Range("S7").Value = 0.00930497082570228
$ws.Range("T7").Value = 0.009304970825702282
